$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 419. This shifts the existing rows
# 419-447 down to 420-448 (matching the diff: new_row[n] == old_row[n-1]).
$ws.Rows.Item(419).Insert()

# Populate the newly inserted row 419 with the new weekly data record.
$ws.Range("A419").Value = 3
$ws.Range("B419").Value = "Femacal de La Calera"
$ws.Range("C419").Value = "Coquimbo"
$ws.Range("D419").Value = 44931
$ws.Range("E419").Value = 5
$ws.Range("F419").Value = 100112012
$ws.Range("G419").Value = "Espinaca"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 120
$ws.Range("K419").Value = 4000
$ws.Range("L419").Value = 4000
$ws.Range("M419").Value = 4000
$ws.Range("N419").Value = "$/docena de atados (3 kilos)"
$ws.Range("O419").Value = "Provincia de Quillota"
$ws.Range("P419").Value = 1333
$ws.Range("Q419").Value = 3
$ws.Range("R419").Value = "Hortaliza"

# Apply the same date-number-format style used by the other Fecha cells.
$ws.Range("D419").NumberFormat = "YYYY-MM-DD HH:MM:SS"
